$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Bmp4"
$ws.Cells.Item(2,3).Value = "Bmpr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 8.675694999999999
$ws.Cells.Item(2,8).Value = 26.027085
$ws.Cells.Item(2,9).Value = 0.5592117158070719
$ws.Cells.Item(2,10).Value = 0.5592117158070719
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 44.547044
$ws.Cells.Item(2,14).Value = 133.641132
$ws.Cells.Item(2,15).Value = 0.3085128263790582
$ws.Cells.Item(2,16).Value = 0.3085128263790582
$ws.Cells.Item(2,17).Value = 386.47656689558
$ws.Cells.Item(2,18).Value = 3478.28910206022
$ws.Cells.Item(2,19).Value = 0.1725239869879224
$ws.Cells.Item(2,20).Value = 0.1725239869879224

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Bmp4"
$ws.Cells.Item(3,3).Value = "Bmpr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 8.675694999999999
$ws.Cells.Item(3,8).Value = 26.027085
$ws.Cells.Item(3,9).Value = 0.5592117158070719
$ws.Cells.Item(3,10).Value = 0.5592117158070719
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 40.23702866666667
$ws.Cells.Item(3,14).Value = 120.711086
$ws.Cells.Item(3,15).Value = 0.2786635952555802
$ws.Cells.Item(3,16).Value = 0.2786635952555802
$ws.Cells.Item(3,17).Value = 349.0841884182566
$ws.Cells.Item(3,18).Value = 3141.75769576431
$ws.Cells.Item(3,19).Value = 0.1558319472358404
$ws.Cells.Item(3,20).Value = 0.1558319472358404

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Bmp4"
$ws.Cells.Item(4,3).Value = "Bmpr2"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 8.675694999999999
$ws.Cells.Item(4,8).Value = 26.027085
$ws.Cells.Item(4,9).Value = 0.5592117158070719
$ws.Cells.Item(4,10).Value = 0.5592117158070719
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 20.15970133333333
$ws.Cells.Item(4,14).Value = 60.47910400000001
$ws.Cells.Item(4,15).Value = 0.13961704029799
$ws.Cells.Item(4,16).Value = 0.13961704029799
$ws.Cells.Item(4,17).Value = 174.8994200590933
$ws.Cells.Item(4,18).Value = 1574.09478053184
$ws.Cells.Item(4,19).Value = 0.07807548466094409
$ws.Cells.Item(4,20).Value = 0.07807548466094411

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Bmp4"
$ws.Cells.Item(5,3).Value = "Bmpr2"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 8.675694999999999
$ws.Cells.Item(5,8).Value = 26.027085
$ws.Cells.Item(5,9).Value = 0.5592117158070719
$ws.Cells.Item(5,10).Value = 0.5592117158070719
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 23.39142333333333
$ws.Cells.Item(5,14).Value = 70.17426999999999
$ws.Cells.Item(5,15).Value = 0.1619984959180618
$ws.Cells.Item(5,16).Value = 0.1619984959180618
$ws.Cells.Item(5,17).Value = 202.9368544558833
$ws.Cells.Item(5,18).Value = 1826.43169010295
$ws.Cells.Item(5,19).Value = 0.09059145686050424
$ws.Cells.Item(5,20).Value = 0.09059145686050424

$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Bmp4"
$ws.Cells.Item(6,3).Value = "Bmpr2"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 8.675694999999999
$ws.Cells.Item(6,8).Value = 26.027085
$ws.Cells.Item(6,9).Value = 0.5592117158070719
$ws.Cells.Item(6,10).Value = 0.5592117158070719
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 16.05764533333333
$ws.Cells.Item(6,14).Value = 48.172936
$ws.Cells.Item(6,15).Value = 0.1112080421493099
$ws.Cells.Item(6,16).Value = 0.1112080421493099
$ws.Cells.Item(6,17).Value = 139.3112333301733
$ws.Cells.Item(6,18).Value = 1253.80109997156
$ws.Cells.Item(6,19).Value = 0.06218884006186073
$ws.Cells.Item(6,20).Value = 0.06218884006186073

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Bmp4"
$ws.Cells.Item(7,3).Value = "Bmpr2"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 6.316050666666666
$ws.Cells.Item(7,8).Value = 18.948152
$ws.Cells.Item(7,9).Value = 0.407115456505913
$ws.Cells.Item(7,10).Value = 0.407115456505913
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 44.547044
$ws.Cells.Item(7,14).Value = 133.641132
$ws.Cells.Item(7,15).Value = 0.3085128263790582
$ws.Cells.Item(7,16).Value = 0.3085128263790582
$ws.Cells.Item(7,17).Value = 281.3613869542293
$ws.Cells.Item(7,18).Value = 2532.252482588064
$ws.Cells.Item(7,19).Value = 0.1256003401492397
$ws.Cells.Item(7,20).Value = 0.1256003401492397

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Bmp4"
$ws.Cells.Item(8,3).Value = "Bmpr2"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 6.316050666666666
$ws.Cells.Item(8,8).Value = 18.948152
$ws.Cells.Item(8,9).Value = 0.407115456505913
$ws.Cells.Item(8,10).Value = 0.407115456505913
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 40.23702866666667
$ws.Cells.Item(8,14).Value = 120.711086
$ws.Cells.Item(8,15).Value = 0.2786635952555802
$ws.Cells.Item(8,16).Value = 0.2786635952555802
$ws.Cells.Item(8,17).Value = 254.1391117347858
$ws.Cells.Item(8,18).Value = 2287.252005613072
$ws.Cells.Item(8,19).Value = 0.1134482567940545
$ws.Cells.Item(8,20).Value = 0.1134482567940545

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Bmp4"
$ws.Cells.Item(9,3).Value = "Bmpr2"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 6.316050666666666
$ws.Cells.Item(9,8).Value = 18.948152
$ws.Cells.Item(9,9).Value = 0.407115456505913
$ws.Cells.Item(9,10).Value = 0.407115456505913
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 20.15970133333333
$ws.Cells.Item(9,14).Value = 60.47910400000001
$ws.Cells.Item(9,15).Value = 0.13961704029799
$ws.Cells.Item(9,16).Value = 0.13961704029799
$ws.Cells.Item(9,17).Value = 127.3296950462009
$ws.Cells.Item(9,18).Value = 1145.967255415808
$ws.Cells.Item(9,19).Value = 0.05684025509692065
$ws.Cells.Item(9,20).Value = 0.05684025509692067

$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Bmp4"
$ws.Cells.Item(10,3).Value = "Bmpr2"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 6.316050666666666
$ws.Cells.Item(10,8).Value = 18.948152
$ws.Cells.Item(10,9).Value = 0.407115456505913
$ws.Cells.Item(10,10).Value = 0.407115456505913
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 23.39142333333333
$ws.Cells.Item(10,14).Value = 70.17426999999999
$ws.Cells.Item(10,15).Value = 0.1619984959180618
$ws.Cells.Item(10,16).Value = 0.1619984959180618
$ws.Cells.Item(10,17).Value = 147.7414149387822
$ws.Cells.Item(10,18).Value = 1329.67273444904
$ws.Cells.Item(10,19).Value = 0.06595209161895299
$ws.Cells.Item(10,20).Value = 0.06595209161895299

$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Bmp4"
$ws.Cells.Item(11,3).Value = "Bmpr2"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 6.316050666666666
$ws.Cells.Item(11,8).Value = 18.948152
$ws.Cells.Item(11,9).Value = 0.407115456505913
$ws.Cells.Item(11,10).Value = 0.407115456505913
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 16.05764533333333
$ws.Cells.Item(11,14).Value = 48.172936
$ws.Cells.Item(11,15).Value = 0.1112080421493099
$ws.Cells.Item(11,16).Value = 0.1112080421493099
$ws.Cells.Item(11,17).Value = 101.4209015126969
$ws.Cells.Item(11,18).Value = 912.788113614272
$ws.Cells.Item(11,19).Value = 0.04527451284674509
$ws.Cells.Item(11,20).Value = 0.04527451284674509

$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Bmp4"
$ws.Cells.Item(12,3).Value = "Bmpr2"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.5224053333333334
$ws.Cells.Item(12,8).Value = 1.567216
$ws.Cells.Item(12,9).Value = 0.03367282768701513
$ws.Cells.Item(12,10).Value = 0.03367282768701512
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 44.547044
$ws.Cells.Item(12,14).Value = 133.641132
$ws.Cells.Item(12,15).Value = 0.3085128263790582
$ws.Cells.Item(12,16).Value = 0.3085128263790582
$ws.Cells.Item(12,17).Value = 23.27161336983467
$ws.Cells.Item(12,18).Value = 209.444520328512
$ws.Cells.Item(12,19).Value = 0.01038849924189604
$ws.Cells.Item(12,20).Value = 0.01038849924189604

$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Bmp4"
$ws.Cells.Item(13,3).Value = "Bmpr2"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.5224053333333334
$ws.Cells.Item(13,8).Value = 1.567216
$ws.Cells.Item(13,9).Value = 0.03367282768701513
$ws.Cells.Item(13,10).Value = 0.03367282768701512
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 40.23702866666667
$ws.Cells.Item(13,14).Value = 120.711086
$ws.Cells.Item(13,15).Value = 0.2786635952555802
$ws.Cells.Item(13,16).Value = 0.2786635952555802
$ws.Cells.Item(13,17).Value = 21.02003837295289
$ws.Cells.Item(13,18).Value = 189.180345356576
$ws.Cells.Item(13,19).Value = 0.009383391225685277
$ws.Cells.Item(13,20).Value = 0.009383391225685277

$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Bmp4"
$ws.Cells.Item(14,3).Value = "Bmpr2"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 0.6666666666666666
$ws.Cells.Item(14,7).Value = 0.5224053333333334
$ws.Cells.Item(14,8).Value = 1.567216
$ws.Cells.Item(14,9).Value = 0.03367282768701513
$ws.Cells.Item(14,10).Value = 0.03367282768701512
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 20.15970133333333
$ws.Cells.Item(14,14).Value = 60.47910400000001
$ws.Cells.Item(14,15).Value = 0.13961704029799
$ws.Cells.Item(14,16).Value = 0.13961704029799
$ws.Cells.Item(14,17).Value = 10.53153549494045
$ws.Cells.Item(14,18).Value = 94.78381945446402
$ws.Cells.Item(14,19).Value = 0.004701300540125265
$ws.Cells.Item(14,20).Value = 0.004701300540125265

$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Bmp4"
$ws.Cells.Item(15,3).Value = "Bmpr2"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 0.6666666666666666
$ws.Cells.Item(15,7).Value = 0.5224053333333334
$ws.Cells.Item(15,8).Value = 1.567216
$ws.Cells.Item(15,9).Value = 0.03367282768701513
$ws.Cells.Item(15,10).Value = 0.03367282768701512
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 23.39142333333333
$ws.Cells.Item(15,14).Value = 70.17426999999999
$ws.Cells.Item(15,15).Value = 0.1619984959180618
$ws.Cells.Item(15,16).Value = 0.1619984959180618
$ws.Cells.Item(15,17).Value = 12.21980430359111
$ws.Cells.Item(15,18).Value = 109.97823873232
$ws.Cells.Item(15,19).Value = 0.005454947438604517
$ws.Cells.Item(15,20).Value = 0.005454947438604516

$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Bmp4"
$ws.Cells.Item(16,3).Value = "Bmpr2"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 0.6666666666666666
$ws.Cells.Item(16,7).Value = 0.5224053333333334
$ws.Cells.Item(16,8).Value = 1.567216
$ws.Cells.Item(16,9).Value = 0.03367282768701513
$ws.Cells.Item(16,10).Value = 0.03367282768701512
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 16.05764533333333
$ws.Cells.Item(16,14).Value = 48.172936
$ws.Cells.Item(16,15).Value = 0.1112080421493099
$ws.Cells.Item(16,16).Value = 0.1112080421493099
$ws.Cells.Item(16,17).Value = 8.388599562908444
$ws.Cells.Item(16,18).Value = 75.49739606617601
$ws.Cells.Item(16,19).Value = 0.003744689240704026
$ws.Cells.Item(16,20).Value = 0.003744689240704026

